# Updated symbol list on Sun Jan  1 15:45:04 UTC 2023 with GitHub Actions
# Applies the refreshed coinranking.com price/volume snapshot to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text LOOKS like a number/percentage -------------------
# These must keep their original `inlineStr`/text nature (not become a
# numeric cell), so each value is entered with a leading apostrophe, which
# is the standard Excel "treat as text" quote-prefix entry.
$numericLookingText = @{
    "D2" = "'244.41"
    "E2" = "'-1.13%"
    "D3" = "'27.46"
    "E3" = "'4.14%"
    "D4" = "'5.055"
    "E4" = "'-0.58%"
    "D5" = "'0.05683"
    "E5" = "'1.16%"
    "E6" = "'-0.72%"
    "D7" = "'0.8213"
    "E7" = "'0.79%"
    "D8" = "'0.8383"
    "E8" = "'-1.26%"
    "D9" = "'0.01001"
    "E9" = "'1,570.90%"
    "D10" = "'0.1325"
    "E10" = "'-1.62%"
    "D11" = "'0.06923"
    "E11" = "'-1.00%"
    "D12" = "'0.02864"
    "E12" = "'1.62%"
    "D13" = "'0.09401"
    "E13" = "'-0.01%"
    "D14" = "'0.001522"
    "E14" = "'0.41%"
    "D15" = "'0.04117"
    "E15" = "'-11.90%"
    "D16" = "'0.006134"
    "E16" = "'-0.30%"
    "D17" = "'3.509"
    "E17" = "'-2.18%"
    "D18" = "'3.001"
    "E18" = "'-1.92%"
    "D19" = "'2.307"
    "E19" = "'8.93%"
    "D20" = "'0.3113"
    "E20" = "'-2.14%"
    "D21" = "'0.03160"
    "E21" = "'-1.02%"
    "D22" = "'0.1292"
    "E22" = "'-2.17%"
    "D23" = "'3.566"
    "E23" = "'-4.71%"
    "D24" = "'0.1374"
    "E24" = "'1.76%"
    "D25" = "'0.001218"
    "E25" = "'-2.61%"
    "D26" = "'0.003870"
    "E26" = "'-16.10%"
    "D27" = "'0.00009799"
    "E27" = "'2.06%"
    "E28" = "'-25.79%"
    "D40" = "'0.03725"
    "E40" = "'1.37%"
    "D41" = "'0.006139"
    "E41" = "'80.74%"
    "D42" = "'0.1052"
    "E42" = "'-23.02%"
    "D43" = "'0.002300"
    "E43" = "'-13.55%"
    "D44" = "'0.009690"
    "E44" = "'12.88%"
    "D45" = "'0.00005216"
    "E45" = "'-1.42%"
    "E46" = "'-0.02%"
    "D48" = "'0.002571"
    "E48" = "'12.90%"
    "E49" = "'-0.02%"
    "E50" = "'-0.02%"
}

# --- Cells whose new text is plain (coin names / links) --------------------
# These don't look like numbers, so a plain assignment keeps them text.
$plainText = @{
    "B9" = "One"
    "C9" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B15" = "CoinExToken"
    "C15" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "B18" = "GateToken"
    "C18" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "B19" = "BTSEToken"
    "C19" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B20" = "BitpandaEcosystemToken"
    "C20" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "B21" = "LiechtensteinCryptoassetsExchange"
    "C21" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "B22" = "ProBitToken"
    "C22" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "B23" = "MCDex"
    "C23" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "B24" = "ZBToken"
    "C24" = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
}

foreach ($addr in $numericLookingText.Keys) {
    $ws.Range($addr).Value = $numericLookingText[$addr]
}

foreach ($addr in $plainText.Keys) {
    $ws.Range($addr).Value = $plainText[$addr]
}
